$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.255456566810608
$ws.Range("B1").Value = 1.915404319763184
$ws.Range("C1").Value = 2.570384502410889
$ws.Range("D1").Value = 3.960100173950195
$ws.Range("E1").Value = 1.119397759437561
